$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Section header
$ws.Range("N1").Value = "Map Construction"

# N column labels (entered top-to-bottom first, skipping the "Name" row)
$ws.Range("N3").Value = "Format"
$ws.Range("N5").Value = "Width"
$ws.Range("N6").Value = "Height"
$ws.Range("N7").Value = "Start X Position"
$ws.Range("N8").Value = "Start Y Position"
$ws.Range("N9").Value = "NPC Count"
$ws.Range("N10").Value = "Map data"

# Data rows for Width/Height/Start X/Start Y/NPC Count
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 50

$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 50

$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 49

$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 49

$ws.Range("O9").Value = 0
$ws.Range("P9").Formula = "=(50*50)-1"

# Map data row n/a values
$ws.Range("O10").Value = "n/a"
$ws.Range("P10").Value = "n/a"

# Header row of the min/max table
$ws.Range("O3").Value = "Min"
$ws.Range("P3").Value = "Max"

# Name row (filled in last)
$ws.Range("N4").Value = "Name (alpha-numeric)"
$ws.Range("O4").Value = "n/a"
$ws.Range("P4").Value = "n/a"

# Apply the same left-aligned style used by column D/I/K to the O and P columns
$ws.Range("O3:P10").HorizontalAlignment = -4131

# Column widths to match target
$ws.Columns.Item(14).ColumnWidth = 21.42578125
$ws.Columns.Item(15).ColumnWidth = 9.140625
$ws.Columns.Item(16).ColumnWidth = 9.140625

# Update the sheet view to match target selection/scroll position
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("M8").Select()
